$wb = $excel.ActiveWorkbook

# Sheet "展览" - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1049
$ws1.Range("F5").Value = 2873
$ws1.Range("F7").Value = 262
$ws1.Range("F11").Value = 124
$ws1.Range("F12").Value = 41
$ws1.Range("F13").Value = 2704
$ws1.Range("F14").Value = 942

# Sheet "全部类型" - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1049
$ws4.Range("F6").Value = 2873
$ws4.Range("F8").Value = 262
$ws4.Range("F13").Value = 124
$ws4.Range("F14").Value = 41
$ws4.Range("F15").Value = 2704
$ws4.Range("F16").Value = 942
